$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 42: Relecture de la documentation
$ws.Range("A42").Value = 44692
$ws.Range("B42").Value = 0.33333333333333331
$ws.Range("C42").Value = 0.36805555555555558
$ws.Range("E42").Value = "Relecture de la documentation"
$ws.Range("F42").Value = "vérification de l'orhographe et de la grammaire`nréimportation d'images qui étaient de mauvaises qualitées"

# Row 43: System de Preselection fonctionnelle
$ws.Range("A43").Value = 44692
$ws.Range("B43").Value = 0.36805555555555558
$ws.Range("C43").Value = 0.45
$ws.Range("E43").Value = "System de Preselection fonctionnelle"

# Row 44: Documentation Selon Model
$ws.Range("A44").Value = 44692
$ws.Range("B44").Value = 0.45
$ws.Range("E44").Value = "Documentation Selon Model"

# Update the sheet view to match the commit (scrolled position / selection)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 39
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E45").Select()
